$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 341.72726
$ws.Range("I6").Value = 257.375
$ws.Range("J6").Value = 566.6667
$ws.Range("K6").Value = 772.125
$ws.Range("L6").Value = 1700.0001
$ws.Range("M6").Value = -660.125
$ws.Range("N6").Value = -1924.0001
$ws.Range("H17").Value = 331.48
$ws.Range("J17").Value = 331.48
$ws.Range("L17").Value = 994.4400000000001
$ws.Range("N17").Value = -1330.44
$ws.Range("H33").Value = 244
$ws.Range("I33").Value = 189.16667
$ws.Range("K33").Value = 189.16667
$ws.Range("M33").Value = 39.83332999999999
$ws.Range("H39").Value = 377.13635
$ws.Range("I39").Value = 119.3
$ws.Range("J39").Value = 592
$ws.Range("K39").Value = 357.9
$ws.Range("L39").Value = 1776
$ws.Range("M39").Value = -61.89999999999998
$ws.Range("N39").Value = -2368
$ws.Range("H43").Value = 17858062
$ws.Range("I43").Value = 71429250
$ws.Range("K43").Value = 71429250
$ws.Range("M43").Value = -71429181
$ws.Range("H70").Value = 1642.0834
$ws.Range("I70").Value = 1790.5
$ws.Range("J70").Value = 900
$ws.Range("K70").Value = 5371.5
$ws.Range("L70").Value = 2700
$ws.Range("M70").Value = -5101.5
$ws.Range("N70").Value = -3240
$ws.Range("H73").Value = 1642.0834
$ws.Range("I73").Value = 1790.5
$ws.Range("J73").Value = 900
$ws.Range("K73").Value = 5371.5
$ws.Range("L73").Value = 2700
$ws.Range("M73").Value = -4435.5
$ws.Range("N73").Value = -4572
$ws.Range("H112").Value = 993.65625
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 993.55554
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 2980.66662
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -5196.66662
$ws.Range("H125").Value = 2560
$ws.Range("J125").Value = 2575
$ws.Range("L125").Value = 23175
$ws.Range("N125").Value = -28095
$ws.Range("H135").Value = 9616589
$ws.Range("I135").Value = 12500939
$ws.Range("J135").Value = 2089.1667
$ws.Range("K135").Value = 112508451
$ws.Range("L135").Value = 18802.5003
$ws.Range("M135").Value = -112505916
$ws.Range("N135").Value = -23872.5003
$ws.Range("H137").Value = 1537.5518
$ws.Range("I137").Value = 1462.421
$ws.Range("J137").Value = 1680.3
$ws.Range("K137").Value = 4387.263
$ws.Range("L137").Value = 5040.9
$ws.Range("M137").Value = -1837.263
$ws.Range("N137").Value = -10140.9
$ws.Range("H138").Value = 2679.375
$ws.Range("I138").Value = 1676
$ws.Range("J138").Value = 3162.4814
$ws.Range("K138").Value = 5028
$ws.Range("L138").Value = 9487.4442
$ws.Range("M138").Value = 112
$ws.Range("N138").Value = -19767.4442

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 23342.857
$ws.Range("J34").Value = 24680
$ws.Range("L34").Value = 24680
$ws.Range("N34").Value = -25222
$ws.Range("H51").Value = 24890
$ws.Range("J51").Value = 24890
$ws.Range("L51").Value = 24890
$ws.Range("N51").Value = -26402
$ws.Range("H132").Value = 1435955.4
$ws.Range("I132").Value = 1028.6061
$ws.Range("J132").Value = 7355028
$ws.Range("K132").Value = 3085.8183
$ws.Range("L132").Value = 22065084
$ws.Range("M132").Value = -555.8182999999999
$ws.Range("N132").Value = -22070144

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 19609558
$ws.Range("I132").Value = 1236.7142
$ws.Range("J132").Value = 33335384
$ws.Range("K132").Value = 3710.1426
$ws.Range("L132").Value = 100006152
$ws.Range("M132").Value = -1180.1426
$ws.Range("N132").Value = -100011212
$ws.Range("H134").Value = 1156.1034
$ws.Range("I134").Value = 999.9545000000001
$ws.Range("J134").Value = 1646.8572
$ws.Range("K134").Value = 2999.8635
$ws.Range("L134").Value = 4940.571599999999
$ws.Range("M134").Value = -464.8635000000004
$ws.Range("N134").Value = -10010.5716

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 35357044
$ws.Range("I5").Value = 44872130
$ws.Range("J5").Value = 15297.143
$ws.Range("K5").Value = 134616390
$ws.Range("L5").Value = 45891.429
$ws.Range("M5").Value = -134616278
$ws.Range("N5").Value = -46115.429
$ws.Range("H114").Value = 577.5357
$ws.Range("I114").Value = 227
$ws.Range("J114").Value = 717.75
$ws.Range("K114").Value = 681
$ws.Range("L114").Value = 2153.25
$ws.Range("M114").Value = 2573
$ws.Range("N114").Value = -8661.25
$ws.Range("H121").Value = 749.7857
$ws.Range("I121").Value = 399.66666
$ws.Range("J121").Value = 845.2727
$ws.Range("K121").Value = 1198.99998
$ws.Range("L121").Value = 2535.8181
$ws.Range("M121").Value = 111.0000199999999
$ws.Range("N121").Value = -5155.8181
$ws.Range("H122").Value = 13591892
$ws.Range("I122").Value = 69444776
$ws.Range("J122").Value = 6056.7026
$ws.Range("K122").Value = 625002984
$ws.Range("L122").Value = 54510.32339999999
$ws.Range("M122").Value = -625000534
$ws.Range("N122").Value = -59410.32339999999
$ws.Range("H126").Value = 83335576
$ws.Range("I126").Value = 166667820
$ws.Range("K126").Value = 500003460
$ws.Range("M126").Value = -499998520
$ws.Range("H130").Value = 2166.6667
$ws.Range("I130").Value = 3000
$ws.Range("K130").Value = 9000
$ws.Range("M130").Value = -3980
$ws.Range("H131").Value = 821.09
$ws.Range("J131").Value = 821.2143
$ws.Range("L131").Value = 2463.6429
$ws.Range("N131").Value = -12543.6429
$ws.Range("H134").Value = 22728082
$ws.Range("I134").Value = 23810132
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 71430396
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -71425326
$ws.Range("N134").Value = -25140
$ws.Range("H135").Value = 35357044
$ws.Range("I135").Value = 44872130
$ws.Range("J135").Value = 15297.143
$ws.Range("K135").Value = 403849170
$ws.Range("L135").Value = 137674.287
$ws.Range("M135").Value = -403846635
$ws.Range("N135").Value = -142744.287

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 15644.211
$ws.Range("J57").Value = 15644.211
$ws.Range("L57").Value = 15644.211
$ws.Range("N57").Value = -17284.211
$ws.Range("H102").Value = 2250.158
$ws.Range("I102").Value = 1575.9
$ws.Range("J102").Value = 2999.3333
$ws.Range("K102").Value = 1575.9
$ws.Range("L102").Value = 2999.3333
$ws.Range("M102").Value = 46.09999999999991
$ws.Range("N102").Value = -6243.3333
$ws.Range("H126").Value = 2282.8462
$ws.Range("I126").Value = 1429
$ws.Range("J126").Value = 2539
$ws.Range("K126").Value = 4287
$ws.Range("L126").Value = 7617
$ws.Range("M126").Value = -1817
$ws.Range("N126").Value = -12557
$ws.Range("H132").Value = 6332.5654
$ws.Range("I132").Value = 1736.1666
$ws.Range("J132").Value = 22879.6
$ws.Range("K132").Value = 5208.4998
$ws.Range("L132").Value = 68638.79999999999
$ws.Range("M132").Value = -2678.4998
$ws.Range("N132").Value = -73698.79999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2536173.5
$ws.Range("I22").Value = 3167191.8
$ws.Range("J22").Value = 12100
$ws.Range("K22").Value = 3167191.8
$ws.Range("L22").Value = 12100
$ws.Range("M22").Value = -3166896.8
$ws.Range("N22").Value = -12690
$ws.Range("H27").Value = 2536173.5
$ws.Range("I27").Value = 3167191.8
$ws.Range("J27").Value = 12100
$ws.Range("K27").Value = 3167191.8
$ws.Range("L27").Value = 12100
$ws.Range("M27").Value = -3167084.8
$ws.Range("N27").Value = -12314
$ws.Range("H40").Value = 13891987
$ws.Range("I40").Value = 3143.2727
$ws.Range("J40").Value = 35717310
$ws.Range("K40").Value = 3143.2727
$ws.Range("L40").Value = 35717310
$ws.Range("M40").Value = -3007.2727
$ws.Range("N40").Value = -35717582
$ws.Range("H46").Value = 2977860.8
$ws.Range("I46").Value = 5209271
$ws.Range("J46").Value = 2647
$ws.Range("K46").Value = 5209271
$ws.Range("L46").Value = 2647
$ws.Range("M46").Value = -5209083
$ws.Range("N46").Value = -3023
$ws.Range("H122").Value = 10300
$ws.Range("I122").Value = 14925
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 44775
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -42325
$ws.Range("N122").Value = -13600
$ws.Range("H132").Value = 43966560
$ws.Range("I132").Value = 87913350
$ws.Range("J132").Value = 19768.23
$ws.Range("K132").Value = 263740050
$ws.Range("L132").Value = 59304.69
$ws.Range("M132").Value = -263737520
$ws.Range("N132").Value = -64364.69
$ws.Range("H136").Value = 82420860
$ws.Range("I136").Value = 67229860
$ws.Range("J136").Value = 111115000
$ws.Range("K136").Value = 201689580
$ws.Range("L136").Value = 333345000
$ws.Range("M136").Value = -201687030
$ws.Range("N136").Value = -333350100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2962.875
$ws.Range("I126").Value = 934.6667
$ws.Range("J126").Value = 4179.8
$ws.Range("K126").Value = 2804.0001
$ws.Range("L126").Value = 12539.4
$ws.Range("M126").Value = -334.0001000000002
$ws.Range("N126").Value = -17479.4
$ws.Range("H132").Value = 30714.553
$ws.Range("I132").Value = 64191.25
$ws.Range("J132").Value = 6367.864
$ws.Range("K132").Value = 192573.75
$ws.Range("L132").Value = 19103.592
$ws.Range("M132").Value = -190043.75
$ws.Range("N132").Value = -24163.592
